# First expenses post, draft
# Adds a "ben_or_claire" column (E) marking a few rows as Claire's, and
# appends four more expense rows (car lockout / rei gear reimbursements).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header + a couple of existing rows tagged in the new E column.
$ws.Range("E1").Value = "ben_or_claire"
$ws.Range("E57").Value = "claire"
$ws.Range("E105").Value = "claire"

# Widen column B so the longer item text (e.g. "air mattress reimbursement")
# is readable.
$ws.Columns.Item(2).ColumnWidth = 25.6640625

# Four new expense rows appended at the bottom of the table.
$ws.Range("A205").Value = 43582
$ws.Range("A205").NumberFormat = "d-mmm"
$ws.Range("B205").Value = "car lockout"
$ws.Range("C205").Value = 65
$ws.Range("D205").Value = "dumb"

$ws.Range("A206").Value = 43510
$ws.Range("A206").NumberFormat = "d-mmm"
$ws.Range("B206").Value = "rei"
$ws.Range("C206").Value = 352.4
$ws.Range("D206").Value = "gear"

$ws.Range("A207").Value = 43512
$ws.Range("A207").NumberFormat = "d-mmm"
$ws.Range("B207").Value = "rei (parents paid)"
$ws.Range("C207").Value = 297.55
$ws.Range("D207").Value = "gear"

$ws.Range("A208").Value = 43613
$ws.Range("A208").NumberFormat = "d-mmm"
$ws.Range("B208").Value = "air mattress reimbursement"
$ws.Range("C208").Value = -151.58
$ws.Range("D208").Value = "gear"

# Match the author's final selection/scroll position.
$ws.Range("E106").Select()
